$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.027.84"
$ws.Range("E2").Value = "  -1.79%  "
$ws.Range("D3").Value = "2.261.50"
$ws.Range("E3").Value = "  -1.75%  "
$ws.Range("E4").Value = "  +0.28%  "
$c = $ws.Range("D5")
$c.Value = "'110.79"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.63%  "
$c = $ws.Range("D6")
$c.Value = "'263.06"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.82%  "
$c = $ws.Range("D7")
$c.Value = "'0.615"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -1.66%  "
$ws.Range("E8").Value = "  +0.08%  "
$c = $ws.Range("D9")
$c.Value = "'0.595"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -3.81%  "
$c = $ws.Range("D10")
$c.Value = "'47.30"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.31%  "
$c = $ws.Range("D11")
$c.Value = "'0.0921"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.04%  "
$c = $ws.Range("D12")
$c.Value = "'8.62"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("E13").Value = "  -0.02%  "
$c = $ws.Range("D14")
$c.Value = "'15.26"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.96%  "
$ws.Range("D15").Value = "2.602.48"
$ws.Range("E15").Value = "  -1.69%  "
$c = $ws.Range("D16")
$c.Value = "'0.846"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("D17").Value = "2.260.93"
$ws.Range("E17").Value = "  -1.82%  "
$ws.Range("D18").Value = "42.943.38"
$ws.Range("E18").Value = "  -1.76%  "
$ws.Range("E19").Value = "  -3.56%  "
$c = $ws.Range("D20")
$c.Value = "'6.77"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +6.86%  "
$c = $ws.Range("D21")
$c.Value = "'70.97"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.Value = "'2.43"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -3.47%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Range("D23")
$c.Value = "'230.12"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.53%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D24")
$c.Value = "'9.54"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.22%  "
$c = $ws.Range("D25")
$c.Value = "'2.83"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.84%  "
$ws.Range("E26").Value = "  +0.02%  "
$c = $ws.Range("D27")
$c.Value = "'11.14"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.29%  "
$c = $ws.Range("D28")
$c.Value = "'3.91"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.02%  "
$c = $ws.Range("D29")
$c.Value = "'40.11"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -6.81%  "
$c = $ws.Range("D30")
$c.Value = "'3.34"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.54%  "
$c = $ws.Range("D31")
$c.Value = "'2.24"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.09%  "
$c = $ws.Range("D32")
$c.Value = "'170.36"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -3.94%  "
$c = $ws.Range("D33")
$c.Value = "'21.03"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.33%  "
$c = $ws.Range("D34")
$c.Value = "'0.0897"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -2.59%  "
$c = $ws.Range("D35")
$c.Value = "'5.67"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("E36").Value = "  -0.88%  "
$c = $ws.Range("D37")
$c.Value = "'4.60"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -4.46%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D38")
$c.Value = "'0.0347"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.47%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D39")
$c.Value = "'3.77"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -6.68%  "
$c = $ws.Range("D40")
$c.Value = "'0.103"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -7.77%  "
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D41")
$c.Value = "'2.47"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +4.15%  "
$c = $ws.Range("D42")
$c.Value = "'74.80"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +7.09%  "
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c = $ws.Range("D43")
$c.Value = "'13.93"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +14.74%  "
$c = $ws.Range("D44")
$c.Value = "'0.231"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -2.42%  "
$c = $ws.Range("D45")
$c.Value = "'6.06"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +11.73%  "
$ws.Range("E46").Value = "  -0.13%  "
$c = $ws.Range("D47")
$c.Value = "'1.35"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.86%  "
$c = $ws.Range("D48")
$c.Value = "'8.53"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -2.32%  "
$c = $ws.Range("D49")
$c.Value = "'0.0984"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -3.08%  "
$c = $ws.Range("D50")
$c.Value = "'99.40"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.18%  "
$c = $ws.Range("D51")
$c.Value = "'1.22"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.28%  "
